# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data block (row 159),
# pushing all existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 159; existing rows 159:185 shift to 160:186
$ws.Rows("159:159").Insert()

# Populate the new row 159 with the new weekly price observation
$ws.Range("A159").Value = 7
$ws.Range("B159").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C159").Value = "Ñuble"
$ws.Range("D159").Value = 44505
$ws.Range("E159").Value = 16
$ws.Range("F159").Value = 100112023
$ws.Range("G159").Value = "Brócoli"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 400
$ws.Range("K159").Value = 650
$ws.Range("L159").Value = 700
$ws.Range("M159").Value = 675
$ws.Range("N159").Value = "$/unidad"
$ws.Range("O159").Value = "Región del Maule"
$ws.Range("P159").Value = 675
$ws.Range("Q159").Value = 1
$ws.Range("R159").Value = "Hortaliza"
